# Updated cryptos list values (price + 1h volume change) per target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.640.10"
$ws.Range("E2").Value = "  -0.31%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.398.87"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.78%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.46"
$ws.Range("E5").Value = "  -0.92%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.96"
$ws.Range("E6").Value = "  +0.75%  "

# Row 7
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  +1.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.404.93"
$ws.Range("E9").Value = "  +0.15%  "

# Row 10
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  -0.57%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  +1.55%  "

# Row 13
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.09"
$ws.Range("E14").Value = "  +0.46%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.836.21"
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.550.25"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.405.09"
$ws.Range("E18").Value = "  +0.26%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.01"
$ws.Range("E19").Value = "  +8.11%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("E20").Value = "  -0.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.11"
$ws.Range("E21").Value = "  +0.43%  "

# Row 22
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.07"
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("E25").Value = "  -2.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.04"
$ws.Range("E26").Value = "  +0.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "567.11"
$ws.Range("E27").Value = "  -1.80%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -4.12%  "

# Row 29
$ws.Range("E29").Value = "  -0.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  +1.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.09"
$ws.Range("E31").Value = "  +2.22%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.33"
$ws.Range("E32").Value = "  -0.86%  "

# Row 33
$ws.Range("E33").Value = "  -1.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.131"
$ws.Range("E34").Value = "  -0.67%  "

# Row 35
$ws.Range("E35").Value = "  -0.58%  "

# Row 36
$ws.Range("E36").Value = "  +3.93%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.23"
$ws.Range("E37").Value = "  +1.78%  "

# Row 38
$ws.Range("E38").Value = "  +1.06%  "

# Row 39
$ws.Range("E39").Value = "  -1.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.28"
$ws.Range("E40").Value = "  +0.06%  "

# Row 41
$ws.Range("E41").Value = "  +0.01%  "

# Row 42
$ws.Range("E42").Value = "  -0.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.66"
$ws.Range("E43").Value = "  +1.29%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  +0.78%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +6.54%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.53"
$ws.Range("E47").Value = "  +0.32%  "

# Row 48
$ws.Range("E48").Value = "  +0.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.590"
$ws.Range("E49").Value = "  +0.25%  "

# Row 50
$ws.Range("E50").Value = "  +0.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.32"
$ws.Range("E51").Value = "  -0.97%  "

